$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New column H ("2022 Traverse") mirrors column G for every row that has
#    a value there (same pattern already used for the other vehicle columns).
# ---------------------------------------------------------------------------
$ws.Range("H1").Value = "2022 Traverse"

$mirrorRows = @(2,5,6,8,9,10,11,14,15,16,17,18,19,20,27,31,39,40,42,45,46,47,48,49,50,54,56)
foreach ($r in $mirrorRows) {
    $ws.Range("H$r").Value = $ws.Range("G$r").Value()
}

# ---------------------------------------------------------------------------
# 2. Harness pin re-labeling for the PT CAN expansion (rows 33-38, col X2).
#    X2-2 / X2-3 : was a plain "X" pin, now labeled "Object CAN +/-"
#    X2-4 / X2-5 : was "CAN +/-", now becomes a plain "X" pin
#    X2-6 / X2-7 : was a plain "X" pin, now labeled "PT CAN +/-"
# ---------------------------------------------------------------------------

# X2-4 -> plain "X" pin (array formula in B, literal X across C:H)
$ws.Range("B35").ClearFormats()
$ws.Range("C35:H35").Value = "X"
$ws.Range("B35").FormulaArray = '=IF(OR(IF(C35:Z35="X",1,0)),"X","")'

# X2-5 -> plain "X" pin (array formula in B, literal X across C:H)
$ws.Range("B36").ClearFormats()
$ws.Range("C36:H36").Value = "X"
$ws.Range("B36").FormulaArray = '=IF(OR(IF(C36:Z36="X",1,0)),"X","")'

# X2-6 -> "PT CAN +"
$ws.Range("B37:H37").Value = "PT CAN +"
$ws.Range("B37").Interior.Color = 15773696

# X2-7 -> "PT CAN -"
$ws.Range("B38:H38").Value = "PT CAN -"
$ws.Range("B38").Interior.Color = 15773696

# X2-2 -> "Object CAN +"
$ws.Range("B33:H33").Value = "Object CAN +"
$ws.Range("B33").Interior.Color = 15773696

# X2-3 -> "Object CAN -"
$ws.Range("B34:H34").Value = "Object CAN -"
$ws.Range("B34").Interior.Color = 15773696

# ---------------------------------------------------------------------------
# 3. Highlight the three "Combined" (D-column) pins affected by the harness
#    change in yellow.
# ---------------------------------------------------------------------------
$ws.Range("D11").Interior.Color = 65535
$ws.Range("D20").Interior.Color = 65535
$ws.Range("D42").Interior.Color = 65535

# ---------------------------------------------------------------------------
# 4. Column widths (best-fit, as Excel recalculates them after the new data
#    is entered).
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 5.83203125
$ws.Columns("B").ColumnWidth = 12
$ws.Columns("C").ColumnWidth = 12
$ws.Columns("E").ColumnWidth = 12
$ws.Columns("G").ColumnWidth = 12
$ws.Columns("H").ColumnWidth = 12.33203125

# ---------------------------------------------------------------------------
# 5. Final selection, matching the author's last cursor position.
# ---------------------------------------------------------------------------
$ws.Range("G37").Select()
